# Update cryptos list (price & 1h volume change) per scraper refresh run.
# D-column price strings use dots as thousands separators and must stay as
# literal text (not be reinterpreted as numbers), so they are entered with a
# leading apostrophe and the cell style is reset to Normal afterwards so no
# stray text-number-format style is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.884.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "'3.440.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'575.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'159.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "'3.441.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "'0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.63%  "
$ws.Range("D10").Value = "'7.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").Value = "'0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'4.035.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'27.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000186"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.77%  "
$ws.Range("D17").Value = "'64.890.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "'3.444.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "'6.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "'13.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "'381.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'7.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.96%  "
$ws.Range("D23").Value = "'0.549"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "'72.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "'0.0000118"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.04%  "
$ws.Range("D27").Value = "'9.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'1.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("D33").Value = "'23.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "'7.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'160.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "'1.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").Value = "'2.914.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "'0.0749"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").Value = "'6.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("D41").Value = "'26.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.76%  "
$ws.Range("D42").Value = "'4.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").Value = "'43.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "'0.0317"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.778"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'26.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("D47").Value = "'2.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").Value = "'317.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'6.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("E51").Value = "  -4.20%  "
